$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove specific cell values (delete content entirely, like the source cells were removed)
$ws.Range("G1").ClearContents()
$ws.Range("E3").ClearContents()
$ws.Range("G4").ClearContents()

# Update the active selection to match the new state
$ws.Range("F11").Select()
